# Get IMACLIM World region out of aggregation table.
#
# The "regions" sheet carries a "World" column (O) that shouldn't be part
# of the aggregation mapping. Clear it out (header + all 49 data rows)
# while leaving every other column/value/style untouched - this is what
# drops the "World" shared string and shifts every other shared-string
# index used across the workbook down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("regions")

[void]$ws.Range("O1:O50").ClearContents()

# The author ended the session on the "regions" sheet, scrolled back to
# the top of the frozen table with the cursor on N1 (USA) - make it the
# active sheet/selection.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N1").Select() | Out-Null
